$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells keep their text formatting so values such as "1.00" or
# "0.0971" are not silently coerced into numbers by Excel.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.707.75"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +2.22%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.044.55"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +3.32%  "

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.10%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "258.86"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +5.45%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.626"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.72%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.06%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "57.80"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -5.24%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.389"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.47%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "57.34"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.12%  "

# Row 11
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.43%  "

# Row 12
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.16%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.83"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.38%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.342.74"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +3.21%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.825"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -2.69%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.46"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -2.68%  "

# Row 17
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.14%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.043.47"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +2.97%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "37.526.43"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +2.01%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "70.27"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.14%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0858"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.56%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.26"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.86%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "229.62"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.25%  "

# Row 24
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +7.61%  "

# Row 25
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.01%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.36"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.68%  "

# Row 27
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.50%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "163.78"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.01%  "

# Row 29
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -5.32%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "20.04"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +2.67%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.36"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.41%  "

# Row 32
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.43%  "

# Row 33
$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0668"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +7.75%  "

# Row 34
$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.77"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.80%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.53"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +11.04%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.54"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.45%  "

# Row 37
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +3.64%  "

# Row 38
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.01%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.82"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +2.46%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.39"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.90%  "

# Row 41
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +4.21%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0971"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.32%  "

# Row 43
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +3.10%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.20"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +1.80%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.413.71"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +3.32%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "16.20"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.25%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "91.67"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.87%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.06"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.40%  "

# Row 49
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +2.98%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.89"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +2.60%  "

# Row 51
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +2.41%  "
